# BOT; UPDATE DATA
#
# Daily-data-refresh bot commit: appends one new day's row (date serial
# 43970 = 2020-05-19) to the three data sheets ("all", "kobe", "other"),
# pushing the trailing footnote row down by one, and moves the "active"
# sheet/view from "all" to "kobe".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": footer row 42 -> 43, new data row 42
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate() | Out-Null

# Push the existing footnote row (currently row 42) down to row 43 and
# open up a fresh row 42 that inherits formatting from the row above it.
$wsAll.Rows.Item(42).Insert(-4121) | Out-Null

$wsAll.Range("A42").Value = 43970
$wsAll.Range("B42").Value = ""
$wsAll.Range("C42").Value = 281
$wsAll.Range("D42").Value = 52
$wsAll.Range("E42").Value = 45
$wsAll.Range("F42").Value = 7
$wsAll.Range("G42").Value = 11
$wsAll.Range("H42").Value = 218

$wsAll.Range("I40").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "kobe": footer row 97 -> 98, new data row 97
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate() | Out-Null

$wsKobe.Rows.Item(97).Insert(-4121) | Out-Null

$wsKobe.Range("A97").Value = 43970
# Daily test-count (column B) isn't published yet for this day, so the
# cell is left fully blank (no style carried over either).
$wsKobe.Range("B97").Clear() | Out-Null
$wsKobe.Range("C97").Value = 2896
$wsKobe.Range("D97").Value = 0
$wsKobe.Range("E97").Value = 283
$wsKobe.Range("F97").Value = 47
$wsKobe.Range("G97").Value = 41
$wsKobe.Range("H97").Value = 6
$wsKobe.Range("I97").Value = 11
$wsKobe.Range("J97").Value = 209

$wsKobe.Range("I75").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "other": footer row 72 -> 73, new data row 72
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate() | Out-Null

$wsOther.Rows.Item(72).Insert(-4121) | Out-Null

$wsOther.Range("A72").Value = 43970
$wsOther.Range("B72").Value = 0
$wsOther.Range("C72").Value = 14
$wsOther.Range("D72").Value = 5
$wsOther.Range("E72").Value = 4
$wsOther.Range("F72").Value = 1
$wsOther.Range("G72").Value = 0
$wsOther.Range("H72").Value = 9

$wsOther.Range("E76").Select() | Out-Null

# ---------------------------------------------------------------------
# Leave "kobe" as the active/selected tab (workbook activeTab=1).
# ---------------------------------------------------------------------
$wsKobe.Activate() | Out-Null
